$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$origStyle = $ws.Range("A1").Style

$ws.Range("D2").Value = '34.609.03'
$ws.Range("E2").Value = '  +1.16%  '
$ws.Range("D3").Value = '1.818.45'
$ws.Range("E3").Value = '  +1.80%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'228.33"
$ws.Range("D5").Style = $origStyle
$ws.Range("E6").Value = '  +1.10%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = "'34.82"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = '  +8.12%  '
$ws.Range("E9").Value = '  +2.20%  '
$ws.Range("E10").Value = '  +1.23%  '
$ws.Range("D12").Value = '2.080.64'
$ws.Range("E12").Value = '  +1.76%  '
$ws.Range("D13").Value = "'11.37"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  +3.45%  '
$ws.Range("D14").Value = '1.823.94'
$ws.Range("E14").Value = '  +2.42%  '
$ws.Range("D15").Value = "'0.648"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  +3.55%  '
$ws.Range("D16").Value = '34.625.86'
$ws.Range("E16").Value = '  +1.26%  '
$ws.Range("D17").Value = "'4.35"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = '  +3.99%  '
$ws.Range("D18").Value = "'69.20"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  +1.84%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.0₃0802'
$ws.Range("E19").Value = '  +0.17%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = "'247.29"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = '  +0.63%  '
$ws.Range("E21").Value = '  +5.85%  '
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("E23").Value = '  +1.53%  '
$ws.Range("D24").Value = "'172.96"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  +6.97%  '
$ws.Range("D25").Value = "'2.08"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  +1.34%  '
$ws.Range("E26").Value = '  +4.24%  '
$ws.Range("D27").Value = "'16.78"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  +2.84%  '
$ws.Range("E28").Value = '  +1.66%  '
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("E30").Value = '  +7.86%  '
$ws.Range("D31").Value = "'0.0532"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  +2.23%  '
$ws.Range("E32").Value = '  +2.95%  '
$ws.Range("E33").Value = '  +1.61%  '
$ws.Range("E34").Value = '  +3.06%  '
$ws.Range("D35").Value = '1.418.76'
$ws.Range("E35").Value = '  -1.34%  '
$ws.Range("E36").Value = '  -1.20%  '
$ws.Range("E37").Value = '  +2.41%  '
$ws.Range("E38").Value = '  +1.74%  '
$ws.Range("E39").Value = '  +1.26%  '
$ws.Range("D40").Value = "'86.23"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  +5.70%  '
$ws.Range("E41").Value = '  +4.48%  '
$ws.Range("D42").Value = "'0.956"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  +3.89%  '
$ws.Range("E43").Value = '  +0.73%  '
$ws.Range("D44").Value = "'13.82"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  -1.99%  '
$ws.Range("E45").Value = '  +1.20%  '
$ws.Range("E46").Value = '  +2.55%  '
$ws.Range("E47").Value = '  +0.84%  '
$ws.Range("D48").Value = '1.981.35'
$ws.Range("E48").Value = '  +2.06%  '
$ws.Range("D49").Value = "'105.89"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  +0.41%  '
$ws.Range("E50").Value = '  +1.69%  '
